$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (no explicit style index / General format) used to restore
# each cell after forcing a temporary text NumberFormat, so numeric-looking
# strings (e.g. "42.636.85", "1.00") are stored as text without leaving the
# cell permanently tagged with a new style index.
$plainStyle = $ws.Range("B2").Style

function Set-TextValue {
    param($CellRef, $Text)
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value2 = $Text
    $rng.Style = $plainStyle
}

Set-TextValue 'D2' '42.636.85'
Set-TextValue 'E2' '  -0.09%  '
Set-TextValue 'D3' '2.527.70'
Set-TextValue 'E3' '  -1.16%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.11%  '
Set-TextValue 'D5' '314.46'
Set-TextValue 'E5' '  +0.63%  '
Set-TextValue 'D6' '98.91'
Set-TextValue 'E6' '  -2.91%  '
Set-TextValue 'D7' '0.563'
Set-TextValue 'E8' '  +0.04%  '
Set-TextValue 'E9' '  -2.52%  '
Set-TextValue 'D10' '35.23'
Set-TextValue 'E10' '  -3.04%  '
Set-TextValue 'D11' '0.0800'
Set-TextValue 'E11' '  -0.50%  '
Set-TextValue 'E12' '  +1.01%  '
Set-TextValue 'E13' '  -2.49%  '
Set-TextValue 'D14' '2.916.63'
Set-TextValue 'E14' '  -1.05%  '
Set-TextValue 'B15' 'WrappedEther'
Set-TextValue 'C15' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D15' '2.526.31'
Set-TextValue 'E15' '  -2.47%  '
Set-TextValue 'B16' 'Chainlink'
Set-TextValue 'C16' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D16' '15.24'
Set-TextValue 'E16' '  -6.02%  '
Set-TextValue 'D17' '0.811'
Set-TextValue 'E17' '  -3.88%  '
Set-TextValue 'D18' '42.659.13'
Set-TextValue 'E18' '  -0.08%  '
Set-TextValue 'D19' '6.61'
Set-TextValue 'E19' '  -3.20%  '
Set-TextValue 'D20' '12.19'
Set-TextValue 'E20' '  -1.60%  '
Set-TextValue 'D21' '0.0₃0940'
Set-TextValue 'E21' '  -1.77%  '
Set-TextValue 'D22' '69.05'
Set-TextValue 'E22' '  -0.16%  '
Set-TextValue 'D23' '241.92'
Set-TextValue 'E23' '  -0.57%  '
Set-TextValue 'E24' '  -1.80%  '
Set-TextValue 'E25' '  -3.53%  '
Set-TextValue 'E26' '  +0.23%  '
Set-TextValue 'D27' '25.47'
Set-TextValue 'E27' '  -3.82%  '
Set-TextValue 'D28' '2.25'
Set-TextValue 'E28' '  -4.63%  '
Set-TextValue 'E29' '  -1.27%  '
Set-TextValue 'E30' '  -6.29%  '
Set-TextValue 'D31' '5.92'
Set-TextValue 'E31' '  +3.85%  '
Set-TextValue 'D32' '156.09'
Set-TextValue 'E32' '  -1.56%  '
Set-TextValue 'E33' '  -2.33%  '
Set-TextValue 'D34' '0.0783'
Set-TextValue 'E34' '  -2.68%  '
Set-TextValue 'D35' '2.64'
Set-TextValue 'E35' '  +0.60%  '
Set-TextValue 'D36' '3.13'
Set-TextValue 'E36' '  -1.98%  '
Set-TextValue 'E37' '  -4.86%  '
Set-TextValue 'D38' '17.58'
Set-TextValue 'E38' '  -3.53%  '
Set-TextValue 'E39' '  -2.52%  '
Set-TextValue 'E40' '  -0.78%  '
Set-TextValue 'E41' '  +0.57%  '
Set-TextValue 'D42' '21.95'
Set-TextValue 'E42' '  -0.66%  '
Set-TextValue 'E43' '  -0.05%  '
Set-TextValue 'B44' 'Maker'
Set-TextValue 'C44' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D44' '2.026.88'
Set-TextValue 'E44' '  +3.10%  '
Set-TextValue 'B45' 'VeChain'
Set-TextValue 'C45' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D45' '0.0296'
Set-TextValue 'E45' '  -1.02%  '
Set-TextValue 'E46' '  -3.76%  '
Set-TextValue 'D47' '8.91'
Set-TextValue 'E47' '  -0.96%  '
Set-TextValue 'D48' '2.769.44'
Set-TextValue 'E48' '  -1.19%  '
Set-TextValue 'D49' '80.05'
Set-TextValue 'E49' '  -1.23%  '
Set-TextValue 'E50' '  -2.70%  '
Set-TextValue 'D51' '71.95'
Set-TextValue 'E51' '  -1.42%  '
